$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ccbdb300564cd7c41f874cad9b627906e54d838b/e2e/a.md"

# --- Overview sheet: update the "Ready for handoff" status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null

# --- zh-cn sheet: record the handback (target file + handback file + handback datetime) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $aMdUrl, "", "", "a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $aMdUrl, "", "", "a.md") | Out-Null

$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-09-01 22:40:37"
$wsZh.Range("K3").Value = "2016-09-01 22:40:37"

$wsZh.Columns.Item(3).AutoFit() | Out-Null
$wsZh.Columns.Item(10).ColumnWidth = $wsZh.Columns.Item(7).ColumnWidth

# --- de-de sheet: record the handback (target file + handback file + handback datetime) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $aMdUrl, "", "", "a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $aMdUrl, "", "", "a.md") | Out-Null

$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDe.Range("K2").Value = "2016-09-01 22:40:45"
$wsDe.Range("K3").Value = "2016-09-01 22:40:45"

$wsDe.Columns.Item(3).AutoFit() | Out-Null
$wsDe.Columns.Item(10).ColumnWidth = $wsDe.Columns.Item(7).ColumnWidth
